# Fruta / hortaliza, semanal
# Insert a new data row at sheet row 275 (shifting the existing rows 275..294
# down to 276..295) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 275..294 down by inserting a fresh row at 275.
$ws.Rows.Item(275).Insert()

# Fill in the new record.
$ws.Cells.Item(275, 1).Value = 10
$ws.Cells.Item(275, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(275, 3).Value = "La Araucanía"
$ws.Cells.Item(275, 4).Value = 44516
$ws.Cells.Item(275, 4).NumberFormat = $ws.Cells.Item(274, 4).NumberFormat
$ws.Cells.Item(275, 5).Value = 9
$ws.Cells.Item(275, 6).Value = "Fruta"
$ws.Cells.Item(275, 7).Value = 100103
$ws.Cells.Item(275, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(275, 9).Value = 100103006
$ws.Cells.Item(275, 10).Value = "Nectarín"
$ws.Cells.Item(275, 11).Value = "Early Glo"
$ws.Cells.Item(275, 12).Value = "Primera"
$ws.Cells.Item(275, 13).Value = 25
$ws.Cells.Item(275, 14).Value = 21000
$ws.Cells.Item(275, 15).Value = 21000
$ws.Cells.Item(275, 16).Value = 21000
$ws.Cells.Item(275, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(275, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(275, 19).Value = 1167
$ws.Cells.Item(275, 20).Value = 18
